$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.697.52"
$ws.Range("E2").Value = "  -4.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.431.32"
$ws.Range("E3").Value = "  -5.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.22"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.39"
$ws.Range("E6").Value = "  -6.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.599"
$ws.Range("E7").Value = "  -4.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.426.23"
$ws.Range("E8").Value = "  -5.13%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("E10").Value = "  -5.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.69"
$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.565"
$ws.Range("E12").Value = "  -8.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "45.91"
$ws.Range("E13").Value = "  -4.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000268"
$ws.Range("E14").Value = "  -4.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.988.25"
$ws.Range("E15").Value = "  -4.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "608.88"
$ws.Range("E16").Value = "  -9.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.21"
$ws.Range("E17").Value = "  -8.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.799.79"
$ws.Range("E18").Value = "  -3.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.434.11"
$ws.Range("E19").Value = "  -5.17%  "

$ws.Range("E20").Value = "  -3.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.04"
$ws.Range("E21").Value = "  -3.96%  "

$ws.Range("E22").Value = "  -4.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.864"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.32"
$ws.Range("E24").Value = "  -9.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "94.68"
$ws.Range("E25").Value = "  -4.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.69"
$ws.Range("E26").Value = "  -5.34%  "

$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.55"
$ws.Range("E28").Value = "  -7.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.86"
$ws.Range("E29").Value = "  -9.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.93"
$ws.Range("E30").Value = "  -7.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("E31").Value = "  -9.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.00"
$ws.Range("E32").Value = "  -8.00%  "

$ws.Range("E33").Value = "  -7.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.71"
$ws.Range("E34").Value = "  -10.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "598.36"
$ws.Range("E35").Value = "  +5.09%  "

$ws.Range("E36").Value = "  -4.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.45"
$ws.Range("E37").Value = "  -3.39%  "

$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  -14.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0995"
$ws.Range("E40").Value = "  -7.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0430"
$ws.Range("E41").Value = "  -4.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.133"
$ws.Range("E42").Value = "  -4.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.343.89"
$ws.Range("E43").Value = "  -5.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.317"
$ws.Range("E44").Value = "  -8.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.06"
$ws.Range("E45").Value = "  -6.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0674"
$ws.Range("E46").Value = "  -7.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("E47").Value = "  -7.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("E48").Value = "  -10.07%  "

$ws.Range("E49").Value = "  -5.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.31"
$ws.Range("E50").Value = "  -3.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.64"
$ws.Range("E51").Value = "  +12.35%  "
